$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -8.570700000000002
$ws.Range("E3").Value = 16.36779999999999
$ws.Range("A12").Value = -21.57770000000001
$ws.Range("D14").Value = -7.768299999999998
$ws.Range("D26").Value = -8.810600000000008
$ws.Range("E30").Value = 15.7173
$ws.Range("D31").Value = -8.888400000000003
$ws.Range("A32").Value = -21.23819999999999
$ws.Range("D35").Value = -8.295500000000002
$ws.Range("A36").Value = -19.919
$ws.Range("D37").Value = -7.838399999999996
$ws.Range("A38").Value = -19.4787
$ws.Range("E44").Value = 16.63759999999998
$ws.Range("D45").Value = -7.651
$ws.Range("A46").Value = -21.7987
$ws.Range("A54").Value = -21.92639999999999
$ws.Range("A55").Value = -22.4434
$ws.Range("D57").Value = -8.530099999999996
$ws.Range("E58").Value = 16.20390000000002
$ws.Range("A67").Value = -21.40179999999998
$ws.Range("A69").Value = -21.59129999999998
$ws.Range("A72").Value = -21.86679999999999
$ws.Range("E84").Value = 16.59419999999999
$ws.Range("E89").Value = 17.34240000000001
$ws.Range("A91").Value = -21.40300000000001
$ws.Range("E91").Value = 17.88690000000003
$ws.Range("E92").Value = 17.98710000000003
$ws.Range("A99").Value = -20.16309999999999
$ws.Range("D100").Value = -8.449199999999998
$ws.Range("D102").Value = -7.767200000000002
$ws.Range("E102").Value = 16.5535
